# divers/backlog.xlsx - "Ajout detection area pour eviter que les balles ne traversent les murs"
#
# The backlog task that used to read "Ajouter des éléments de décor " (row 13,
# column A) is reworded to "Ajouter des éléments de décor destructibles" -
# this is the new backlog entry for the "ball/wall collision" feature
# described in the commit message. Excel re-numbers the shared-string table
# when the old text becomes unused and the new text is appended, which is
# exactly what ripples through every other A/D cell's shared-string index in
# the diff - we only need to touch the content cells, not the indices.
#
# Row 30 ("Ajouter animation de rechargement") gains a DRI ("Qui") and an OK
# ("Fait") value, matching the pattern used by the other completed rows.
#
# Finally the active selection moves from A17 to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Ajouter des éléments de décor destructibles"

$ws.Range("B30").Value = "DRI"
$ws.Range("C30").Value = "OK"

$ws.Range("A4").Select()
